# Sea creature refactoring:
# Add a new "Sheet2" after "Sheet1", populate it with the Blue/Red value
# table (with shared formulas), and make Sheet2 the active/selected sheet
# (mirroring the removal of tabSelected from Sheet1's sheetView).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Header row
$ws2.Range("A1").Value = "Blueaug"
$ws2.Range("C1").Value = "Redval"
$ws2.Range("D1").Value = "Blueval1"
$ws2.Range("E1").Value = "Blueval2"

# Reference values
$ws2.Range("A2").Value = 1
$ws2.Range("A3").Value = 0

# C column inputs
$ws2.Range("C2").Value = 1
$ws2.Range("C3").Value = 2
$ws2.Range("C4").Value = 3
$ws2.Range("C5").Value = 4
$ws2.Range("C6").Value = 5
$ws2.Range("C7").Value = 5
$ws2.Range("C8").Value = 4
$ws2.Range("C9").Value = 3
$ws2.Range("C10").Value = 2
$ws2.Range("C11").Value = 1

# D / E formulas, entered in row order so shared formula groups line up
# the same way they do in the target workbook (si=0,1,2,3).
$ws2.Range("D2").Formula = "=(C2*2)+50 -A`$2"
$ws2.Range("E2").Formula = "=(C2*2)+50 -A`$3"

$ws2.Range("D3:D6").Formula = "=(C3*2)+50 -A`$2"
$ws2.Range("E3:E6").Formula = "=(C3*2)+50 -A`$3"

$ws2.Range("D7").Formula = "=(C7*2)+50 +A`$2"
$ws2.Range("E7").Formula = "=(C7*2)+50 +A`$3"

$ws2.Range("D8:D11").Formula = "=(C8*2)+50 +A`$2"
$ws2.Range("E8:E11").Formula = "=(C8*2)+50 +A`$3"

# Sheet2 becomes the active sheet/tab, with E5 selected.
$ws2.Range("E5").Select()
